# Denmark Division 1 - update league bases (2024-01-31 20:02)
#
# The underlying match-result rows had their per-match B:AC data (match id,
# home/away team, scores, result, odds, PL columns) reassigned between a
# handful of same-kickoff-date rows. Column A (the running match index) and
# the Div/Div-Original-Name columns (C/D, identical across all rows here)
# stay put; every other column from B through AC moves wholesale between the
# listed rows. For rows 308/309/310 this is a 3-way rotation, the rest are
# simple pairwise swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC as 1-based indices (A=1 is left untouched).
$cols = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)

function Get-RowValues($sheet, $row, $colList) {
    $values = @{}
    foreach ($c in $colList) {
        $values[$c] = $sheet.Cells.Item($row, $c).Value2
    }
    return $values
}

function Set-RowValues($sheet, $row, $colList, $values) {
    foreach ($c in $colList) {
        $sheet.Cells.Item($row, $c).Value = $values[$c]
    }
}

# Snapshot every affected row's current (pre-edit) B:AC content before any
# writes happen, so later writes never read already-overwritten data.
$rowsToSnapshot = @(280, 281, 308, 309, 310, 337, 338, 377, 378, 439, 440, 540, 541, 586, 587)
$snapshot = @{}
foreach ($r in $rowsToSnapshot) {
    $snapshot[$r] = Get-RowValues $ws $r $cols
}

# destination row -> source row whose (pre-edit) content it should receive
$moveMap = @{
    280 = 281
    281 = 280
    308 = 310
    309 = 308
    310 = 309
    337 = 338
    338 = 337
    377 = 378
    378 = 377
    439 = 440
    440 = 439
    540 = 541
    541 = 540
    586 = 587
    587 = 586
}

foreach ($dest in $moveMap.Keys) {
    $src = $moveMap[$dest]
    Set-RowValues $ws $dest $cols $snapshot[$src]
}
